$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "path" column (E): these values had a stray leading "//" — they
# should use a single leading "/" (root-relative XPath), for rows 2-9.
$ws.Range("E2").Value = "/ProcedimientoXunta/SI460A_4/Delegacion/cmbDelegacion"
$ws.Range("E3").Value = "/ProcedimientoXunta/SI460A_4/Delegacion/txtIdioma"
$ws.Range("E4").Value = "/ProcedimientoXunta/SI460A/Notificacion/rblModalidad"
$ws.Range("E5").Value = "/ProcedimientoXunta/SI460A/Notificacion/rbNotificar"
$ws.Range("E6").Value = "/ProcedimientoXunta/CT102A/Director/txtNombre"
$ws.Range("E7").Value = "/ProcedimientoXunta/CT102A/Director/txtApel1"
$ws.Range("E8").Value = "/ProcedimientoXunta/CT102A/Director/txtApel2"
$ws.Range("E9").Value = "/ProcedimientoXunta/CT102A/Director/txtCp"

# Add two new rows describing new fields that distinguish radio buttons
# (RB) and checkboxes (Checks) from plain text fields, as per the commit
# message "Distingue entre RB, Checks y campos de texto + formateoXML".
$ws.Range("A10").Value = "sexoDirector"
$ws.Range("B10").Value = "Sexo director"
$ws.Range("C10").Value = "Sexo director"
$ws.Range("D10").Value = "Integer"
$ws.Range("E10").Value = "/ProcedimientoXunta/CT102A/Director/rbSexo"
$ws.Range("F10").Value = "ESP-sexDirector"

$ws.Range("A11").Value = "sabeIngles"
$ws.Range("B11").Value = "Sabe inglés"
$ws.Range("C11").Value = "Sabe inglés"
$ws.Range("D11").Value = "Boolean"
$ws.Range("E11").Value = "/ProcedimientoXunta/CT102A/Director/chSabe"
$ws.Range("F11").Value = "ESP-englishKnowledge"

# Move the active selection to A12, right below the new data (matches the
# author's final cursor position in the saved workbook).
$ws.Range("A12").Select()
